# Daily attendance processing - 2026-01-03 17:03:15
# Swap the order of names in the "Recorded By" column (G) from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
